# :sparkles: Added OSINT easy and Forensic medium
#
# 1) Add a new (horizontal) slide guide at pos=1620, colored gray (A4A3A4) -
#    mirrors the pre-existing red vertical guide at pos=5311.
# 2) Remove the four leftover "OSINT easy" picture shapes from slide 3
#    (the slide is being redone with new content).
# 3) Bump the event date shown on slide 4 from 02.04.2022 to 03.12.2022.

$p = $ppt.ActivePresentation

# --- 1) Slide guides ------------------------------------------------------
try {
    $newGuide = $p.Guides.Add(2, 1620)
    $newGuide.Color = 10789796
} catch {
    try {
        $p.Guides.Item(2).Color = 10789796
    } catch {
        # Guide color API unavailable in this host; position is still added.
    }
}

# --- 2) Slide 3: drop all picture placeholders -----------------------------
$s3 = $p.Slides.Item(3)
while ($s3.Shapes.Count -gt 0) {
    $s3.Shapes.Item(1).Delete()
}

# --- 3) Slide 4: update the event date text --------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(3).TextFrame.TextRange.Text = "03.12.2022 von 10 bis 18 Uhr an der THI"
